$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Das vs APCSP"
$ws.Range("B2").Value = "Das wins"
$ws.Range("C2").Value = "Criminal aka Fraud"
$ws.Range("D2").Value = "Criminal"
$ws.Range("E2").Value = "Adk"
$ws.Range("F2").Value = "Newmans OFfice"
